$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 18:22"

# --- Turquia overtakes China in ranking order (row 10 / row 11 swap labels) ---
$ws.Range("A10").Value = "Turquia"
$ws.Range("A11").Value = "China"

# --- Moldavia overtakes Banglades in ranking order (row 59 / row 60 swap labels) ---
$ws.Range("A59").Value = "Moldavia"
$ws.Range("A60").Value = "Banglades"

# --- Updated statistics for Estados Unidos (row 4) ---
$ws.Range("B4").Value = 741999
$ws.Range("C4").Value = 3207
$ws.Range("D4").Value = 68658
$ws.Range("E4").Value = 633700
$ws.Range("G4").Value = 627
$ws.Range("H4").Value = 39641

# --- Updated statistics for Italia (row 6) ---
$ws.Range("B6").Value = 178972
$ws.Range("C6").Value = 3047
$ws.Range("D6").Value = 47055
$ws.Range("E6").Value = 108257
$ws.Range("F6").Value = 2635
$ws.Range("G6").Value = 433
$ws.Range("H6").Value = 23660

# --- Updated statistics for Francia (row 7) ---
$ws.Range("B7").Value = 152578
$ws.Range("C7").Value = 785
$ws.Range("E7").Value = 96877
$ws.Range("G7").Value = 395
$ws.Range("H7").Value = 19718

# --- Updated statistics for row 10 (now Turquia) ---
$ws.Range("B10").Value = 86306
$ws.Range("C10").Value = 3977
$ws.Range("D10").Value = 11976
$ws.Range("E10").Value = 72313
$ws.Range("F10").Value = 1922
$ws.Range("G10").Value = 127
$ws.Range("H10").Value = 2017

# --- Updated statistics for row 11 (now China) ---
$ws.Range("B11").Value = 82735
$ws.Range("C11").Value = 16
$ws.Range("D11").Value = 77062
$ws.Range("E11").Value = 1041
$ws.Range("F11").Value = 85
$ws.Range("H11").Value = 4632

# --- Updated statistics for India (row 20) ---
$ws.Range("B20").Value = 17137
$ws.Range("C20").Value = 772
$ws.Range("D20").Value = 2769
$ws.Range("E20").Value = 13818

# --- Updated statistics for Rumania (row 32) ---
$ws.Range("E32").Value = 6409
$ws.Range("G32").Value = 24
$ws.Range("H32").Value = 445

# --- Updated statistics for Pakistan (row 33) ---
$ws.Range("B33").Value = 8348
$ws.Range("C33").Value = 710
$ws.Range("E33").Value = 6312
$ws.Range("G33").Value = 25
$ws.Range("H33").Value = 168

# --- Updated statistics for Bielorrusia (row 47) ---
$ws.Range("D47").Value = 494
$ws.Range("E47").Value = 4238

# --- Updated statistics for Argelia (row 58) ---
$ws.Range("B58").Value = 2629
$ws.Range("C58").Value = 95
$ws.Range("D58").Value = 1047
$ws.Range("E58").Value = 1207
$ws.Range("F58").Value = 40
$ws.Range("G58").Value = 8
$ws.Range("H58").Value = 375

# --- Updated statistics for row 59 (now Moldavia) ---
$ws.Range("B59").Value = 2472
$ws.Range("C59").Value = 121
$ws.Range("D59").Value = 457
$ws.Range("E59").Value = 1948
$ws.Range("F59").Value = 80
$ws.Range("G59").Value = 10
$ws.Range("H59").Value = 67

# --- Updated statistics for row 60 (now Banglades) ---
$ws.Range("B60").Value = 2456
$ws.Range("C60").Value = 312
$ws.Range("D60").Value = 75
$ws.Range("E60").Value = 2290
$ws.Range("F60").Value = 1
$ws.Range("G60").Value = 7
$ws.Range("H60").Value = 91

# --- Updated statistics for Gibraltar (row 134) ---
$ws.Range("D134").Value = 120
$ws.Range("E134").Value = 12

# --- Updated statistics for Gabon (row 138) ---
$ws.Range("B138").Value = 109
$ws.Range("C138").Value = 1
$ws.Range("E138").Value = 101

# --- Updated statistics for Siria (row 167) ---
$ws.Range("B167").Value = 39
$ws.Range("C167").Value = 1
$ws.Range("G167").Value = 1
$ws.Range("H167").Value = 3
